$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header renames: ht_goals_h -> HTHG, ht_goals_a -> HTAG
$ws.Range("I1").Value = "HTHG"
$ws.Range("J1").Value = "HTAG"

# Row group swap: [175, 177]
# row 175 <= original row 177 data
$ws.Range("B175").Value = 7302795
$ws.Range("C175").Value = "Peru Liga 1"
$ws.Range("D175").Value = 45221.70833333334
$ws.Range("E175").Value = "Unin Comercio"
$ws.Range("F175").Value = "Deportivo Garcilaso"
$ws.Range("G175").Value = 1
$ws.Range("H175").Value = 2
$ws.Range("I175").Value = 1
$ws.Range("J175").Value = 0
$ws.Range("K175").Value = "A"
$ws.Range("L175").Value = 2.25
$ws.Range("M175").Value = 3.3
$ws.Range("N175").Value = 2.7
$ws.Range("O175").Value = 1.75
$ws.Range("P175").Value = 3.6
$ws.Range("Q175").Value = 4
$ws.Range("R175").Value = -0.5
$ws.Range("S175").Value = 1.8
$ws.Range("T175").Value = 2
$ws.Range("U175").Value = 2.75
$ws.Range("V175").Value = 1.825
$ws.Range("W175").Value = 1.975
$ws.Range("X175").Value = -1
$ws.Range("Y175").Value = -1
$ws.Range("Z175").Value = 3
$ws.Range("AA175").Value = -1
$ws.Range("AB175").Value = 1
$ws.Range("AC175").Value = 0.4125
$ws.Range("AD175").Value = -0.5

# row 177 <= original row 175 data
$ws.Range("B177").Value = 7302200
$ws.Range("C177").Value = "Peru Liga 1"
$ws.Range("D177").Value = 45221.70833333334
$ws.Range("E177").Value = "Carlos Manucci"
$ws.Range("F177").Value = "Deportivo Binacional"
$ws.Range("G177").Value = 3
$ws.Range("H177").Value = 2
$ws.Range("I177").Value = 0
$ws.Range("J177").Value = 0
$ws.Range("K177").Value = "H"
$ws.Range("L177").Value = 2
$ws.Range("M177").Value = 3.2
$ws.Range("N177").Value = 3.75
$ws.Range("O177").Value = 1.75
$ws.Range("P177").Value = 3.4
$ws.Range("Q177").Value = 4.333
$ws.Range("R177").Value = -0.5
$ws.Range("S177").Value = 1.85
$ws.Range("T177").Value = 1.95
$ws.Range("U177").Value = 2.5
$ws.Range("V177").Value = 1.85
$ws.Range("W177").Value = 1.95
$ws.Range("X177").Value = 0.75
$ws.Range("Y177").Value = -1
$ws.Range("Z177").Value = -1
$ws.Range("AA177").Value = 0.8500000000000001
$ws.Range("AB177").Value = -1
$ws.Range("AC177").Value = 0.8500000000000001
$ws.Range("AD177").Value = -1

# Row group swap: [183, 185, 186]
# row 183 <= original row 185 data
$ws.Range("B183").Value = 7384626
$ws.Range("C183").Value = "Peru Liga 1"
$ws.Range("D183").Value = 45228.70833333334
$ws.Range("E183").Value = "Sporting Cristal"
$ws.Range("F183").Value = "Alianza Atletico"
$ws.Range("G183").Value = 3
$ws.Range("H183").Value = 0
$ws.Range("I183").Value = 3
$ws.Range("J183").Value = 0
$ws.Range("K183").Value = "H"
$ws.Range("L183").Value = 1.3
$ws.Range("M183").Value = 5
$ws.Range("N183").Value = 9
$ws.Range("O183").Value = 1.166
$ws.Range("P183").Value = 6.5
$ws.Range("Q183").Value = 13
$ws.Range("R183").Value = -2
$ws.Range("S183").Value = 1.85
$ws.Range("T183").Value = 1.95
$ws.Range("U183").Value = 3.25
$ws.Range("V183").Value = 2
$ws.Range("W183").Value = 1.8
$ws.Range("X183").Value = 0.1659999999999999
$ws.Range("Y183").Value = -1
$ws.Range("Z183").Value = -1
$ws.Range("AA183").Value = 0.8500000000000001
$ws.Range("AB183").Value = -1
$ws.Range("AC183").Value = -0.5
$ws.Range("AD183").Value = 0.4

# row 185 <= original row 186 data
$ws.Range("B185").Value = 7384627
$ws.Range("C185").Value = "Peru Liga 1"
$ws.Range("D185").Value = 45228.70833333334
$ws.Range("E185").Value = "Universitario de Deportes"
$ws.Range("F185").Value = "Sport Huancayo"
$ws.Range("G185").Value = 2
$ws.Range("H185").Value = 0
$ws.Range("I185").Value = 1
$ws.Range("J185").Value = 0
$ws.Range("K185").Value = "H"
$ws.Range("L185").Value = 1.25
$ws.Range("M185").Value = 5
$ws.Range("N185").Value = 12
$ws.Range("O185").Value = 1.181
$ws.Range("P185").Value = 6
$ws.Range("Q185").Value = 13
$ws.Range("R185").Value = -1.75
$ws.Range("S185").Value = 1.8
$ws.Range("T185").Value = 2
$ws.Range("U185").Value = 2.75
$ws.Range("V185").Value = 1.85
$ws.Range("W185").Value = 1.95
$ws.Range("X185").Value = 0.181
$ws.Range("Y185").Value = -1
$ws.Range("Z185").Value = -1
$ws.Range("AA185").Value = 0.4
$ws.Range("AB185").Value = -0.5
$ws.Range("AC185").Value = -1
$ws.Range("AD185").Value = 0.95

# row 186 <= original row 183 data
$ws.Range("B186").Value = 7384628
$ws.Range("C186").Value = "Peru Liga 1"
$ws.Range("D186").Value = 45228.70833333334
$ws.Range("E186").Value = "Deportivo Binacional"
$ws.Range("F186").Value = "FBC Melgar"
$ws.Range("G186").Value = 1
$ws.Range("H186").Value = 2
$ws.Range("I186").Value = 1
$ws.Range("J186").Value = 1
$ws.Range("K186").Value = "A"
$ws.Range("L186").Value = 2.75
$ws.Range("M186").Value = 3.3
$ws.Range("N186").Value = 2.375
$ws.Range("O186").Value = 3.3
$ws.Range("P186").Value = 3.6
$ws.Range("Q186").Value = 2
$ws.Range("R186").Value = 0.5
$ws.Range("S186").Value = 1.8
$ws.Range("T186").Value = 2
$ws.Range("U186").Value = 2.75
$ws.Range("V186").Value = 1.975
$ws.Range("W186").Value = 1.875
$ws.Range("X186").Value = -1
$ws.Range("Y186").Value = -1
$ws.Range("Z186").Value = 1
$ws.Range("AA186").Value = -1
$ws.Range("AB186").Value = 1
$ws.Range("AC186").Value = 0.4875
$ws.Range("AD186").Value = -0.5

# Row group swap: [228, 229]
# row 228 <= original row 229 data
$ws.Range("B228").Value = 7818817
$ws.Range("C228").Value = "Peru Liga 1"
$ws.Range("D228").Value = 45346.70833333334
$ws.Range("E228").Value = "Sport Boys"
$ws.Range("F228").Value = "Cusco FC"
$ws.Range("G228").Value = 3
$ws.Range("H228").Value = 0
$ws.Range("I228").Value = 2
$ws.Range("J228").Value = 0
$ws.Range("K228").Value = "H"
$ws.Range("L228").Value = 2.2
$ws.Range("M228").Value = 3.2
$ws.Range("N228").Value = 3.2
$ws.Range("O228").Value = 1.6
$ws.Range("P228").Value = 3.75
$ws.Range("Q228").Value = 5.75
$ws.Range("R228").Value = -0.75
$ws.Range("S228").Value = 1.85
$ws.Range("T228").Value = 2
$ws.Range("U228").Value = 2.5
$ws.Range("V228").Value = 1.975
$ws.Range("W228").Value = 1.875
$ws.Range("X228").Value = 0.6000000000000001
$ws.Range("Y228").Value = -1
$ws.Range("Z228").Value = -1
$ws.Range("AA228").Value = 0.8500000000000001
$ws.Range("AB228").Value = -1
$ws.Range("AC228").Value = 0.9750000000000001
$ws.Range("AD228").Value = -1

# row 229 <= original row 228 data
$ws.Range("B229").Value = 7818816
$ws.Range("C229").Value = "Peru Liga 1"
$ws.Range("D229").Value = 45346.70833333334
$ws.Range("E229").Value = "UTC Cajamarca"
$ws.Range("F229").Value = "Universitario de Deportes"
$ws.Range("G229").Value = 0
$ws.Range("H229").Value = 0
$ws.Range("I229").Value = 0
$ws.Range("J229").Value = 0
$ws.Range("K229").Value = "D"
$ws.Range("L229").Value = 3.3
$ws.Range("M229").Value = 3.3
$ws.Range("N229").Value = 2.1
$ws.Range("O229").Value = 4.5
$ws.Range("P229").Value = 3.2
$ws.Range("Q229").Value = 1.95
$ws.Range("R229").Value = 0.5
$ws.Range("S229").Value = 2
$ws.Range("T229").Value = 1.85
$ws.Range("U229").Value = 2
$ws.Range("V229").Value = 1.775
$ws.Range("W229").Value = 2.1
$ws.Range("X229").Value = -1
$ws.Range("Y229").Value = 2.2
$ws.Range("Z229").Value = -1
$ws.Range("AA229").Value = 1
$ws.Range("AB229").Value = -1
$ws.Range("AC229").Value = -1
$ws.Range("AD229").Value = 1.1

# Row group swap: [312, 313]
# row 312 <= original row 313 data
$ws.Range("B312").Value = 8086337
$ws.Range("C312").Value = "Peru Liga 1"
$ws.Range("D312").Value = 45416.72916666666
$ws.Range("E312").Value = "Comerciantes Unidos"
$ws.Range("F312").Value = "Union Comercio"
$ws.Range("G312").Value = 3
$ws.Range("H312").Value = 1
$ws.Range("I312").Value = 1
$ws.Range("J312").Value = 0
$ws.Range("K312").Value = "H"
$ws.Range("L312").Value = 1.615
$ws.Range("M312").Value = 3.75
$ws.Range("N312").Value = 5.5
$ws.Range("O312").Value = 2.2
$ws.Range("P312").Value = 3.2
$ws.Range("Q312").Value = 3.3
$ws.Range("R312").Value = -0.25
$ws.Range("S312").Value = 1.925
$ws.Range("T312").Value = 1.875
$ws.Range("U312").Value = 2.25
$ws.Range("V312").Value = 1.875
$ws.Range("W312").Value = 1.925
$ws.Range("X312").Value = 1.2
$ws.Range("Y312").Value = -1
$ws.Range("Z312").Value = -1
$ws.Range("AA312").Value = 0.925
$ws.Range("AB312").Value = -1
$ws.Range("AC312").Value = 0.875
$ws.Range("AD312").Value = -1

# row 313 <= original row 312 data
$ws.Range("B313").Value = 8086253
$ws.Range("C313").Value = "Peru Liga 1"
$ws.Range("D313").Value = 45416.72916666666
$ws.Range("E313").Value = "Carlos Manucci"
$ws.Range("F313").Value = "Deportivo Garcilaso"
$ws.Range("G313").Value = 1
$ws.Range("H313").Value = 1
$ws.Range("I313").Value = 1
$ws.Range("J313").Value = 0
$ws.Range("K313").Value = "D"
$ws.Range("L313").Value = 2
$ws.Range("M313").Value = 3.2
$ws.Range("N313").Value = 3.6
$ws.Range("O313").Value = 2.5
$ws.Range("P313").Value = 3
$ws.Range("Q313").Value = 2.8
$ws.Range("R313").Value = 0
$ws.Range("S313").Value = 1.775
$ws.Range("T313").Value = 2.025
$ws.Range("U313").Value = 2.5
$ws.Range("V313").Value = 2
$ws.Range("W313").Value = 1.8
$ws.Range("X313").Value = -1
$ws.Range("Y313").Value = 2
$ws.Range("Z313").Value = -1
$ws.Range("AA313").Value = 0
$ws.Range("AB313").Value = 0
$ws.Range("AC313").Value = -1
$ws.Range("AD313").Value = 0.8

# Row group swap: [338, 340]
# row 338 <= original row 340 data
$ws.Range("B338").Value = 8240876
$ws.Range("C338").Value = "Peru Liga 1"
$ws.Range("D338").Value = 45437.70833333334
$ws.Range("E338").Value = "Universitario de Deportes"
$ws.Range("F338").Value = "CD Los Chankas"
$ws.Range("G338").Value = 4
$ws.Range("H338").Value = 0
$ws.Range("I338").Value = 2
$ws.Range("J338").Value = 0
$ws.Range("K338").Value = "H"
$ws.Range("L338").Value = 1.1
$ws.Range("M338").Value = 7.5
$ws.Range("N338").Value = 23
$ws.Range("O338").Value = 1.062
$ws.Range("P338").Value = 11
$ws.Range("Q338").Value = 23
$ws.Range("R338").Value = -2.75
$ws.Range("S338").Value = 1.875
$ws.Range("T338").Value = 1.925
$ws.Range("U338").Value = 3.75
$ws.Range("V338").Value = 1.975
$ws.Range("W338").Value = 1.825
$ws.Range("X338").Value = 0.06200000000000006
$ws.Range("Y338").Value = -1
$ws.Range("Z338").Value = -1
$ws.Range("AA338").Value = 0.875
$ws.Range("AB338").Value = -1
$ws.Range("AC338").Value = 0.4875
$ws.Range("AD338").Value = -0.5

# row 340 <= original row 338 data
$ws.Range("B340").Value = 8240870
$ws.Range("C340").Value = "Peru Liga 1"
$ws.Range("D340").Value = 45437.70833333334
$ws.Range("E340").Value = "Deportivo Garcilaso"
$ws.Range("F340").Value = "FBC Melgar"
$ws.Range("G340").Value = 1
$ws.Range("H340").Value = 3
$ws.Range("I340").Value = 0
$ws.Range("J340").Value = 0
$ws.Range("K340").Value = "A"
$ws.Range("L340").Value = 2.625
$ws.Range("M340").Value = 3.2
$ws.Range("N340").Value = 2.75
$ws.Range("O340").Value = 3.1
$ws.Range("P340").Value = 3.2
$ws.Range("Q340").Value = 2.25
$ws.Range("R340").Value = 0.25
$ws.Range("S340").Value = 1.875
$ws.Range("T340").Value = 1.925
$ws.Range("U340").Value = 2.5
$ws.Range("V340").Value = 1.85
$ws.Range("W340").Value = 1.95
$ws.Range("X340").Value = -1
$ws.Range("Y340").Value = -1
$ws.Range("Z340").Value = 1.25
$ws.Range("AA340").Value = -1
$ws.Range("AB340").Value = 0.925
$ws.Range("AC340").Value = 0.8500000000000001
$ws.Range("AD340").Value = -1
